$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank scores for rows 21 and 22 (batch process overdue items)
$ws.Range("C21").Value = 15
$ws.Range("C22").Value = 10

# Update the view: scroll window so row 13 is the top row (topLeftCell A13),
# then set the active selection cell
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F18").Select()
